$d = $word.ActiveDocument

$pairs = @(
    @("73×60=", "46×39="),
    @("73×35=", "21×32="),
    @("34×71=", "33×57="),
    @("13×27=", "28×87="),
    @("31×22=", "32×46="),
    @("71×41=", "84×34="),
    @("96×42=", "33×40="),
    @("20×67=", "47×34="),
    @("56×48=", "25×23="),
    @("25×99=", "33×37="),
    @("97×12=", "56×50="),
    @("30×65=", "79×34="),
    @("60×11=", "18×88="),
    @("33×98=", "45×18="),
    @("18×68=", "21×60="),
    @("87×45=", "30×61="),
    @("16×27=", "43×47="),
    @("83×24=", "90×69="),
    @("39×34=", "87×70="),
    @("97×90=", "19×57="),
    @("65×78=", "79×64="),
    @("30×18=", "19×49="),
    @("33×80=", "12×37="),
    @("35×20=", "41×34="),
    @("94×93=", "61×41=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
